$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-categorize E2/F2 from "dimension" to "measure"
$ws.Range("E2").Value = "iaest-measure:subseccion-descripcion"
$ws.Range("F2").Value = "iaest-measure:sector-vab-descripcion"

# Update row 3 ("dim" -> "medida") for columns E and F, matching D3/G3 semantics
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Update row 4 ("skos:Concept" -> "xsd:int") for columns E and F
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Remove the mapping file references in row 5 for columns E and F
$ws.Range("E5").Clear()
$ws.Range("F5").Clear()
